$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 132 (shifts old rows 132-152 down to 134-154)
$ws.Range("A132:R133").Insert()

# --- New row 132: Apio, Primera, fecha 2021-11-11 (44511) ---
$ws.Cells.Item(132, 1).Value = 9
$ws.Cells.Item(132, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(132, 3).Value = "Metropolitana"
$ws.Cells.Item(132, 4).Value = 44511
$ws.Cells.Item(132, 5).Value = 13
$ws.Cells.Item(132, 6).Value = 100112017
$ws.Cells.Item(132, 7).Value = "Apio"
$ws.Cells.Item(132, 8).Value = "Americana (o)"
$ws.Cells.Item(132, 9).Value = "Primera"
$ws.Cells.Item(132, 10).Value = 79
$ws.Cells.Item(132, 11).Value = 8000
$ws.Cells.Item(132, 12).Value = 9000
$ws.Cells.Item(132, 13).Value = 8494
$ws.Cells.Item(132, 14).Value = "$/docena de matas"
$ws.Cells.Item(132, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(132, 16).Value = 1416
$ws.Cells.Item(132, 17).Value = 6
$ws.Cells.Item(132, 18).Value = "Hortaliza"

# --- New row 133: Apio, Segunda, fecha 2021-11-11 (44511) ---
$ws.Cells.Item(133, 1).Value = 9
$ws.Cells.Item(133, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(133, 3).Value = "Metropolitana"
$ws.Cells.Item(133, 4).Value = 44511
$ws.Cells.Item(133, 5).Value = 13
$ws.Cells.Item(133, 6).Value = 100112017
$ws.Cells.Item(133, 7).Value = "Apio"
$ws.Cells.Item(133, 8).Value = "Americana (o)"
$ws.Cells.Item(133, 9).Value = "Segunda"
$ws.Cells.Item(133, 10).Value = 35
$ws.Cells.Item(133, 11).Value = 6000
$ws.Cells.Item(133, 12).Value = 6000
$ws.Cells.Item(133, 13).Value = 6000
$ws.Cells.Item(133, 14).Value = "$/docena de matas"
$ws.Cells.Item(133, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(133, 16).Value = 1000
$ws.Cells.Item(133, 17).Value = 6
$ws.Cells.Item(133, 18).Value = "Hortaliza"

# Make sure the date cells keep the date/time number format used by column D elsewhere
$ws.Range("D132:D133").NumberFormat = $ws.Range("D131").NumberFormat
